$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a handful of "error" values in column F (rows 1-25 block) ---
# Row 6 (RM 21): F was blank -> now has a value
$ws.Range("F6").Value = 16.43
# Row 8 (RM 38): F had a value -> now blank
$ws.Range("F8").Value = ""
# Row 19 (RM 125): F was blank -> now has a value
$ws.Range("F19").Value = 17.81
# Row 21 (RM 135): F had a value -> now blank
$ws.Range("F21").Value = ""
# Row 23 (RM 140): F was blank -> now has a value
$ws.Range("F23").Value = 16.48

# --- Remove the "RM 232" row entirely (row 26); rows below shift up by one ---
$ws.Range("A26:F26").EntireRow.Delete()

# --- Remove the "SC 92" row entirely (now at row 27 after the previous delete) ---
$ws.Range("A27:F27").EntireRow.Delete()

# --- Fix up remaining value differences in the shifted SC rows ---
# Row 26 is now "SC 5": D column stays, E becomes blank
$ws.Range("E26").Value = ""
# Row 27 is now "SC 101": E gets a value, F becomes blank
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = ""
# Row 29 is now "SC 119": E becomes blank, F gets a value
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = 18.06
